$wb = $excel.ActiveWorkbook

# --- Sheet "instituicoes": append rows 6-9 in column A ---
$wsInst = $wb.Worksheets.Item("instituicoes")
$wsInst.Range("A6").Value = "UNIRIO"
$wsInst.Range("A7").Value = "Liquigás"
$wsInst.Range("A8").Value = "BNDES"
$wsInst.Range("A9").Value = "CEFET-RJ"

# --- Sheet "assuntos": append rows 4-9 in columns A, B, C ---
$wsAssuntos = $wb.Worksheets.Item("assuntos")

$wsAssuntos.Range("A4").Value = "Acentos Diferenciais"
$wsAssuntos.Range("B4").Value = "Português"
$wsAssuntos.Range("C4").Value = 3

$wsAssuntos.Range("A5").Value = "Ortografia Oficial"
$wsAssuntos.Range("B5").Value = "Português"
$wsAssuntos.Range("C5").Value = 4

$wsAssuntos.Range("A6").Value = "Emprego do Hífen"
$wsAssuntos.Range("B6").Value = "Português"
$wsAssuntos.Range("C6").Value = 5

$wsAssuntos.Range("A7").Value = "Expressões Problemáticas"
$wsAssuntos.Range("B7").Value = "Português"
$wsAssuntos.Range("C7").Value = 6

$wsAssuntos.Range("A8").Value = "Substantivo"
$wsAssuntos.Range("B8").Value = "Português"
$wsAssuntos.Range("C8").Value = 7

$wsAssuntos.Range("A9").Value = "Adjetivo"
$wsAssuntos.Range("B9").Value = "Português"
$wsAssuntos.Range("C9").Value = 8
